# Adds a new "case_03" sheet (a copy of case_01's data, with the
# "direction" column switched from forward "f" to backward "b" for the
# two rows that represent the newly-added backward reactions), documents
# it on the README sheet, and leaves the new sheet as the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Document the new test case on the README sheet ------------------
$readme = $wb.Worksheets.Item("README")
$readme.Range("A5").Value = "case_03"
$readme.Range("B5").Value = "same as case 1, but include backward reactions"

# --- 2. Create case_03 from a copy of case_01, placed after case_02 -----
$case01 = $wb.Worksheets.Item("case_01")
$case02 = $wb.Worksheets.Item("case_02")
$case01.Copy($null, $case02)

$case03 = $wb.Worksheets.Item($case02.Index + 1)
$case03.Name = "case_03"

# kcat_b is now taken into account: two reactions also run backward.
$case03.Range("C5").Value = "b"
$case03.Range("C10").Value = "b"

# --- 3. Update selections / active sheet so case_03 is the active tab ---
$readme.Activate()
[void]$readme.Range("I9").Select()

$case03.Activate()
[void]$case03.Range("G19").Select()
